$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('Construction Materials(8)', 0.6668117385209243),
    @('Multiline Retail(7)', 0.637158856603528),
    @('Multi-Utilities(18)', 0.5909373411224138),
    @('Containers & Packaging(12)', 0.5820552138888447),
    @('Machinery(85)', 0.5773156141689123),
    @('Road & Rail(22)', 0.5710044321418385),
    @('Auto Components(21)', 0.5537157241117897),
    @('Air Freight & Logistics(11)', 0.5118791143994464),
    @('Energy Equipment & Services(32)', 0.5006236236380753),
    @('Electric Utilities(28)', 0.5001381885300904),
    @('Airlines(14)', 0.4917499218936111),
    @('Specialty Retail(58)', 0.4879427848495779),
    @('Trading Companies & Distributors(25)', 0.4854130637614995),
    @('Gas Utilities(12)', 0.4786730779818642),
    @('Wireless Telecommunication Services(14)', 0.4702411239861925),
    @('Hotels, Restaurants & Leisure(50)', 0.4685900344946026),
    @('Chemicals(51)', 0.4634824288800502),
    @('Media(42)', 0.4632595642572921),
    @('Construction & Engineering(20)', 0.4551704769440679),
    @('Semiconductors & Semiconductor Equipment(68)', 0.4341880161375474),
    @('Leisure Products(11)', 0.4273113347199632),
    @('Building Products(23)', 0.4193499308950515),
    @('Household Durables(39)', 0.4173637377507899),
    @('Capital Markets(75)', 0.4068347560461073),
    @('Household Products(10)', 0.4064015974987705),
    @('Marine(15)', 0.4033073324332358),
    @('Aerospace & Defense(37)', 0.3975633353773633),
    @('Communications Equipment(45)', 0.3803307251461573),
    @('Oil, Gas & Consumable Fuels(122)', 0.3755636762179344),
    @('Insurance(75)', 0.3678833024531595),
    @('Commercial Services & Supplies(52)', 0.3607899127967236),
    @('Consumer Finance(15)', 0.3521445803575914),
    @('Textiles, Apparel & Luxury Goods(29)', 0.3515114243143465),
    @('Diversified Consumer Services(17)', 0.3499089347939527),
    @('Metals & Mining(89)', 0.3388671756232035),
    @('Water Utilities(12)', 0.3368986889980476),
    @('Diversified Telecommunication Services(20)', 0.3291281176173813),
    @('Professional Services(35)', 0.3241847688660761),
    @('Electrical Equipment(28)', 0.3177397643706062),
    @('Life Sciences Tools & Services(19)', 0.3115688160483179),
    @('Banks(246)', 0.3029498278140831),
    @('Food & Staples Retailing(15)', 0.3026795069584989),
    @('Software(66)', 0.2980981953635595),
    @('Internet & Direct Marketing Retail(15)', 0.2952069581098644),
    @('Health Care Providers & Services(46)', 0.2907406161623396),
    @('IT Services(52)', 0.2895560816946511),
    @('Beverages(21)', 0.270311358287111),
    @('Personal Products(19)', 0.2646998600836321),
    @('Health Care Equipment & Supplies(83)', 0.2530052926919384),
    @('Thrifts & Mortgage Finance(47)', 0.2454875116583609),
    @('Entertainment(22)', 0.2426693261706819),
    @('Food Products(44)', 0.1977162084373913),
    @('Pharmaceuticals(48)', 0.1616154705497828),
    @('Biotechnology(126)', 0.1573268451086477)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = [double]$data[$i][1]
}

# Remove now-unused trailing rows (previously rows 56-59)
$ws.Range("A56:B59").Clear()
